# Apply updated inventory cost results from server run.
$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 56665.56446366477
$ws.Range("O2").Value = 55688.3537912803

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 37870.89546016529
$ws.Range("I2").Value = 226534.2213208755
$ws.Range("L2").Value = 142476.3005577463
$ws.Range("M2").Value = 93731.81239819515
$ws.Range("N2").Value = 27145.45968928043
$ws.Range("O2").Value = 40525.81511985242

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 23418.88254333185
$ws.Range("B2").Value = 17864.57798134877
$ws.Range("E2").Value = 88591.48758750911
$ws.Range("I2").Value = 120309.5553113526
$ws.Range("M2").Value = 29233.1369406852
$ws.Range("N2").Value = 35850.72954647117
$ws.Range("O2").Value = 21420.44672873601

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 832.9248897141631

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 27405.00233090627
$ws.Range("N2").Value = 4217.51601927484
$ws.Range("O2").Value = 18378.03620052791
